# Applies the "Fixed bugs, esp. in Specs." edit to the Reference sheet's
# Param_Table (D2:G8 -> D2:G9): a new "Max_Build / PCT" parameter row is
# inserted after "CO2_MTons" (row 4), which bumps Capital Cost / Fixed Cost /
# Variable Cost / CO2 Cost each down one slot, with a fresh "CO2 Cost" entry
# re-added at the new bottom row (row 9). The two named ranges that describe
# that table are widened to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined names: grow Param_Table / Params by one row ------------------
$wb.Names.Item("Param_Table").RefersTo = "=Reference!`$D`$2:`$G`$9"
$wb.Names.Item("Params").RefersTo = "=Reference!`$D`$2:`$D`$5"

# --- Row 4 (CO2_MTons): index 2 -> 6 ---------------------------------------
$ws.Range("E4").Value = 6

# --- Row 5 becomes the new "Max_Build" entry -------------------------------
$ws.Range("D5").Value = "Max_Build"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "PCT"
$ws.Range("G5").Value = 1

# --- Row 6: was Fixed Cost, now Capital Cost (index 4 -> 2) ----------------
$ws.Range("D6").Value = "Capital Cost"
$ws.Range("E6").Value = 2

# --- Row 7: was Variable Cost, now Fixed Cost (index 5 -> 3) ---------------
$ws.Range("D7").Value = "Fixed Cost"
$ws.Range("E7").Value = 3

# --- Row 8: was CO2 Cost, now Variable Cost (index 6 -> 4) -----------------
$ws.Range("D8").Value = "Variable Cost"
$ws.Range("E8").Value = 4

# --- Row 9 (new): re-add CO2 Cost at the bottom of the table ---------------
$ws.Range("D9").Value = "CO2 Cost"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "B$"
$ws.Range("G9").Formula = "=1/Thousand"

# --- Stray formatted-but-empty cell G10 is no longer present ---------------
$ws.Range("G10").ClearContents()

# --- The empty G-column filler strip now reaches row 20 --------------------
$ws.Range("G19").Copy($ws.Range("G20"))

# --- Selection moved to H9 --------------------------------------------------
$ws.Range("H9").Select()
